$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New demo row showing a heavy (medium) bottom border style, appended
# after the existing "right_border" example row (row 37 -> new row 39).
$cell = $ws.Range("A39")
$cell.Value = "heavy_bottom_border"
$cell.RowHeight = 16

$cell.Borders.Item(9).LineStyle = 1
$cell.Borders.Item(9).Weight = -4138

$ws.Range("A9").Select()
